# Scheduled-runner data refresh: rewrite the computed market-price / profit
# columns (H:N) for the affected Leve rows on each job sheet, matching the
# latest pull from the price feed.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 3333.75
$ws.Range("I98").Value = 3445
$ws.Range("J98").Value = 3000
$ws.Range("K98").Value = 3445
$ws.Range("L98").Value = 3000
$ws.Range("M98").Value = -1947
$ws.Range("N98").Value = -5996
$ws.Range("H113").Value = 4157.091
$ws.Range("I113").Value = 3265
$ws.Range("J113").Value = 4355.3335
$ws.Range("K113").Value = 3265
$ws.Range("L113").Value = 4355.3335
$ws.Range("M113").Value = -11
$ws.Range("N113").Value = -10863.3335
$ws.Range("H122").Value = 3333.75
$ws.Range("I122").Value = 3445
$ws.Range("J122").Value = 3000
$ws.Range("K122").Value = 10335
$ws.Range("L122").Value = 9000
$ws.Range("M122").Value = -7885
$ws.Range("N122").Value = -13900
$ws.Range("H137").Value = 4122.1904
$ws.Range("I137").Value = 4345.0586
$ws.Range("J137").Value = 3175
$ws.Range("K137").Value = 13035.1758
$ws.Range("L137").Value = 9525
$ws.Range("M137").Value = -10485.1758
$ws.Range("N137").Value = -14625

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H28").Value = 16836
$ws.Range("I28").Value = 10980.125
$ws.Range("K28").Value = 10980.125
$ws.Range("M28").Value = -10788.125
$ws.Range("H99").Value = 16836
$ws.Range("I99").Value = 10980.125
$ws.Range("K99").Value = 10980.125
$ws.Range("M99").Value = -7985.125

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H35").Value = 37500
$ws.Range("J35").Value = 37500
$ws.Range("L35").Value = 37500
$ws.Range("N35").Value = -38120
$ws.Range("H98").Value = 0
$ws.Range("J98").Value = 0
$ws.Range("L98").ClearContents()
$ws.Range("N98").Value = 0

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 85067.664
$ws.Range("I99").Value = 201102.4
$ws.Range("J99").Value = 2185.7144
$ws.Range("K99").Value = 201102.4
$ws.Range("L99").Value = 2185.7144
$ws.Range("M99").Value = -199604.4
$ws.Range("N99").Value = -5181.7144
$ws.Range("H126").Value = 85067.664
$ws.Range("I126").Value = 201102.4
$ws.Range("J126").Value = 2185.7144
$ws.Range("K126").Value = 603307.2
$ws.Range("L126").Value = 6557.1432
$ws.Range("M126").Value = -600837.2
$ws.Range("N126").Value = -11497.1432
$ws.Range("H131").Value = 47000
$ws.Range("J131").Value = 47000
$ws.Range("L131").Value = 47000
$ws.Range("N131").Value = -57080

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H38").Value = 38.8
$ws.Range("I38").Value = 26.222221
$ws.Range("J38").Value = 71.14286
$ws.Range("K38").Value = 78.666663
$ws.Range("L38").Value = 213.42858
$ws.Range("M38").Value = 268.333337
$ws.Range("N38").Value = -907.42858
$ws.Range("H110").Value = 2363.4
$ws.Range("I110").Value = 939
$ws.Range("J110").Value = 4500
$ws.Range("K110").Value = 2817
$ws.Range("L110").Value = 13500
$ws.Range("M110").Value = 1273
$ws.Range("N110").Value = -21680
$ws.Range("H113").Value = 14286288
$ws.Range("I113").Value = 26316360
$ws.Range("J113").Value = 578.375
$ws.Range("K113").Value = 78949080
$ws.Range("L113").Value = 1735.125
$ws.Range("M113").Value = -78946910
$ws.Range("N113").Value = -6075.125
$ws.Range("H131").Value = 1191.0819
$ws.Range("J131").Value = 1272.1091
$ws.Range("L131").Value = 3816.3273
$ws.Range("N131").Value = -13896.3273

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H96").Value = 20497.166
$ws.Range("J96").Value = 20497.166
$ws.Range("L96").Value = 20497.166
$ws.Range("N96").Value = -25989.166
$ws.Range("H102").Value = 3302.0588
$ws.Range("I102").Value = 1479.375
$ws.Range("J102").Value = 4922.222
$ws.Range("K102").Value = 1479.375
$ws.Range("L102").Value = 4922.222
$ws.Range("M102").Value = 142.625
$ws.Range("N102").Value = -8166.222
$ws.Range("H122").Value = 1754.2307
$ws.Range("I122").Value = 1690.5
$ws.Range("J122").Value = 1966.6666
$ws.Range("K122").Value = 5071.5
$ws.Range("L122").Value = 5899.9998
$ws.Range("M122").Value = -2621.5
$ws.Range("N122").Value = -10799.9998
$ws.Range("H126").Value = 3109
$ws.Range("I126").Value = 3109
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 9327
$ws.Range("L126").Value = 0
$ws.Range("M126").ClearContents()
$ws.Range("N126").Value = -6857

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3230
$ws.Range("I7").Value = 2983.3333
$ws.Range("J7").Value = 3600
$ws.Range("K7").Value = 2983.3333
$ws.Range("L7").Value = 3600
$ws.Range("M7").Value = -2871.3333
$ws.Range("N7").Value = -3824
$ws.Range("H40").Value = 26512.5
$ws.Range("I40").Value = 26512.5
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 26512.5
$ws.Range("L40").Value = 0
$ws.Range("M40").ClearContents()
$ws.Range("N40").Value = -26376.5
$ws.Range("H96").Value = 5000000
$ws.Range("J96").Value = 5000000
$ws.Range("L96").Value = 5000000
$ws.Range("N96").Value = -5005492
$ws.Range("H126").Value = 3230
$ws.Range("I126").Value = 2983.3333
$ws.Range("J126").Value = 3600
$ws.Range("K126").Value = 8949.999899999999
$ws.Range("L126").Value = 10800
$ws.Range("M126").Value = -6479.999899999999
$ws.Range("N126").Value = -15740
$ws.Range("H132").Value = 12731.315
$ws.Range("I132").Value = 2846.3076
$ws.Range("J132").Value = 34148.832
$ws.Range("K132").Value = 8538.9228
$ws.Range("L132").Value = 102446.496
$ws.Range("M132").Value = -6008.9228
$ws.Range("N132").Value = -107506.496

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 2205.5625
$ws.Range("I81").Value = 1162.6364
$ws.Range("J81").Value = 4500
$ws.Range("K81").Value = 2325.2728
$ws.Range("L81").Value = 9000
$ws.Range("M81").Value = -1264.2728
$ws.Range("N81").Value = -11122
$ws.Range("H84").Value = 2205.5625
$ws.Range("I84").Value = 1162.6364
$ws.Range("J84").Value = 4500
$ws.Range("K84").Value = 11626.364
$ws.Range("L84").Value = 45000
$ws.Range("M84").Value = -6322.364000000001
$ws.Range("N84").Value = -55608
$ws.Range("H107").Value = 552.6667
$ws.Range("I107").Value = 334.66666
$ws.Range("J107").Value = 661.6667
$ws.Range("K107").Value = 1003.99998
$ws.Range("L107").Value = 1985.0001
$ws.Range("M107").Value = 916.0000200000001
$ws.Range("N107").Value = -5825.0001
$ws.Range("H113").Value = 510.25
$ws.Range("I113").Value = 300.33334
$ws.Range("J113").Value = 1140
$ws.Range("K113").Value = 901.0000200000001
$ws.Range("L113").Value = 3420
$ws.Range("M113").Value = 1268.99998
$ws.Range("N113").Value = -7760
$ws.Range("H122").Value = 2437.3704
$ws.Range("I122").Value = 2477.3333
$ws.Range("J122").Value = 2297.5
$ws.Range("K122").Value = 7431.999899999999
$ws.Range("L122").Value = 6892.5
$ws.Range("M122").Value = -4981.999899999999
$ws.Range("N122").Value = -11792.5
$ws.Range("H126").Value = 1089.4642
$ws.Range("I126").Value = 815.25
$ws.Range("J126").Value = 1775
$ws.Range("K126").Value = 2445.75
$ws.Range("L126").Value = 5325
$ws.Range("M126").Value = 24.25
$ws.Range("N126").Value = -10265
